$d = $word.ActiveDocument

$s = $d.Styles.Item("Normal")

# Font: Times New Roman, size 12pt (24 half-points)
$s.Font.Name = "Times New Roman"
$s.Font.Size = 12

# Paragraph: 1.5 line spacing, justified alignment
$s.ParagraphFormat.LineSpacingRule = [Microsoft.Office.Interop.Word.WdLineSpacing]::wdLineSpaceMultiple
$s.ParagraphFormat.LineSpacing = 18
$s.ParagraphFormat.Alignment = [Microsoft.Office.Interop.Word.WdParagraphAlignment]::wdAlignParagraphJustify
